# Update the "Crypto" worksheet with the expanded asset list (25 rows
# instead of 15), refreshed mean-return / risk figures, and the new
# dimension extent A1:C26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Crypto")

$data = @(
    @("ADA-USD",   -0.0009462549077120909, 0.03630537429437645),
    @("ATOM-USD",   0.0004004341971715087, 0.04330960451409918),
    @("AVAX-USD",  -0.0009475469251185598, 0.04485220147448776),
    @("BCH-USD",    0.001407352743383836,  0.04430072809446627),
    @("BNB-USD",   -0.0002666370955085417, 0.02895642791321457),
    @("BTC-USD",    0.0001698017340538461, 0.02728721768065267),
    @("CRO-USD",   -0.001935140215022296,  0.03557390366834624),
    @("DOGE-USD",   0.0005055873766004,    0.04694877812707692),
    @("DOT-USD",   -0.001076934762782215,  0.03583406844817972),
    @("ETC-USD",    0.0005457174036354702, 0.04918708106781185),
    @("ETH-USD",    0.0004779069833942835, 0.03619237636843207),
    @("KAS-USD",    0.01744103067694832,   0.1230996776299637),
    @("LEO-USD",   -0.0002853039963843259, 0.02466974165848193),
    @("LINK-USD",   0.001232268108914391,  0.04220231508622439),
    @("LTC-USD",    0.0009022075568424984, 0.04038740999764472),
    @("MATIC-USD",  0.0009886723621352019, 0.05065994433348981),
    @("OKB-USD",    0.003389285708800178,  0.03770878311670382),
    @("SHIB-USD",   0.0004775359105210174, 0.05382868312718332),
    @("SOL-USD",    0.0002351162794294515, 0.05479576322261236),
    @("TON-USD",   -0.0001231138256954889, 0.04772686211557989),
    @("TRX-USD",    0.0004388946013284949, 0.02490479007018209),
    @("WBTC-USD",   0.000179030553937395,  0.02758758242718384),
    @("XLM-USD",    0.0004373874387499149, 0.04372835900203738),
    @("XMR-USD",   -0.0001819086296563999, 0.02952524701840095),
    @("XRP-USD",    0.001573038748370876,  0.04954929129429662)
)

# Reference cell whose style (bold font, border, centered alignment) should
# be replicated onto every asset-name cell in column A, including the newly
# added rows.
$styleSource = $ws.Cells.Item(2, 1)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $styleSource.Copy($cellA)
    $cellA.Value = $entry[0]

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
